$wb = $excel.ActiveWorkbook

$oldGuid = "8b2a69ec-4bfe-459f-964d-48d51a1ccbbc"
$newGuid = "3c28b869-dcf7-4306-bca0-f83f22c7d64d"
$oldHash = "9c6d62a270e38b8bdc471336282fd1ae54654f88"
$newHash = "cda3d4e51277e48a2b66d894e50de819a42620ee"

$newMd = "$newGuid.md"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

# Sheet "Overview" - md file name and handoff date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = "2016-03-24 21:17:15"
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = $newMd
}

# Sheet "zh-cn" - md file name, xlf file name, handoff datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = $newMd
$wsZhCn.Range("D2").Value = $newZhXlf
$wsZhCn.Range("E2").Value = "2016-03-24 21:17:11"
$i = 0
foreach ($h in $wsZhCn.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) {
        $h.TextToDisplay = $newMd
    } else {
        $h.TextToDisplay = $newZhXlf
    }
}

# Sheet "de-de" - md file name, xlf file name
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = $newMd
$wsDeDe.Range("D2").Value = $newDeXlf
$j = 0
foreach ($h in $wsDeDe.Hyperlinks) {
    $j = $j + 1
    if ($j -eq 1) {
        $h.TextToDisplay = $newMd
    } else {
        $h.TextToDisplay = $newDeXlf
    }
}
